$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update references for row 2 (order 0)
$ws.Range("B2").Value = "C:\Users\IRDGFRM\Downloads\20-04-2022_09h-22m.pdf"

# Force ship_out_date to stay plain text (Excel would otherwise read "08.2022"
# as a number), then restore the cell to its original unstyled state.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "08.2022"
$ws.Range("F2").Style = "Normal"

# Update references for row 3 (order 1)
$ws.Range("B3").Value = "C:\Users\IRDGFRM\Downloads\20-04-2022_09h-22m.pdf"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "16.2022"
$ws.Range("F3").Style = "Normal"

# Remove rows 4, 5 and 6 entirely
$ws.Rows("4:6").Delete()
